$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "LastName" / "FirstName" header values in B1 and C1
$b1 = $ws.Range("B1").Value()
$c1 = $ws.Range("C1").Value()
$ws.Range("B1").Value = $c1
$ws.Range("C1").Value = $b1

# Decrease CNE values in A2:A11 by 10
for ($r = 2; $r -le 11; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cur = $cell.Value()
    $cell.Value = $cur - 10
}

# Update the sheet view: scroll so A1 is the top-left cell, and select E9
$ws.Range("A1").Select()
$ws.Range("E9").Select()
